$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F27").Value = 972.2
$ws.Range("G27").Value = 1735.14

$ws.Range("G28").Select()
